# Auto-generated edit script: updates cached market-board derived values
# across multiple profession sheets (per upstream commit "chore: update Sheets via scheduled runner").
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1257.2106
$ws.Range("J17").Value = 1277.6111
$ws.Range("L17").Value = 3832.8333
$ws.Range("N17").Value = -4168.8333
$ws.Range("H74").Value = 4163.5
$ws.Range("I74").Value = 4396.2
$ws.Range("J74").Value = 3000
$ws.Range("K74").Value = 4396.2
$ws.Range("L74").Value = 3000
$ws.Range("M74").Value = -3460.2
$ws.Range("N74").Value = -4872
$ws.Range("H76").Value = 3930
$ws.Range("I76").Value = 2866.6667
$ws.Range("J76").Value = 4993.3335
$ws.Range("K76").Value = 2866.6667
$ws.Range("L76").Value = 4993.3335
$ws.Range("M76").Value = -2551.6667
$ws.Range("N76").Value = -5623.3335
$ws.Range("H77").Value = 4163.5
$ws.Range("I77").Value = 4396.2
$ws.Range("J77").Value = 3000
$ws.Range("K77").Value = 21981
$ws.Range("L77").Value = 15000
$ws.Range("M77").Value = -17301
$ws.Range("N77").Value = -24360
$ws.Range("H79").Value = 3930
$ws.Range("I79").Value = 2866.6667
$ws.Range("J79").Value = 4993.3335
$ws.Range("K79").Value = 2866.6667
$ws.Range("L79").Value = 4993.3335
$ws.Range("M79").Value = -1774.6667
$ws.Range("N79").Value = -7177.3335
$ws.Range("H112").Value = 3005.625
$ws.Range("J112").Value = 3199.2856
$ws.Range("L112").Value = 9597.856800000001
$ws.Range("N112").Value = -11813.8568

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1642.1578
$ws.Range("I2").Value = 1224.6923
$ws.Range("J2").Value = 2546.6667
$ws.Range("K2").Value = 1224.6923
$ws.Range("L2").Value = 2546.6667
$ws.Range("M2").Value = -1111.6923
$ws.Range("N2").Value = -2772.6667
$ws.Range("H45").Value = 953.3333
$ws.Range("I45").Value = 860
$ws.Range("K45").Value = 860
$ws.Range("M45").Value = -483
$ws.Range("H63").Value = 280067.28
$ws.Range("I63").Value = 314816.56
$ws.Range("J63").Value = 2073
$ws.Range("K63").Value = 314816.56
$ws.Range("L63").Value = 2073
$ws.Range("M63").Value = -314130.56
$ws.Range("N63").Value = -3445
$ws.Range("H66").Value = 280067.28
$ws.Range("I66").Value = 314816.56
$ws.Range("J66").Value = 2073
$ws.Range("K66").Value = 1574082.8
$ws.Range("L66").Value = 10365
$ws.Range("M66").Value = -1570650.8
$ws.Range("N66").Value = -17229
$ws.Range("H116").Value = 1642.1578
$ws.Range("I116").Value = 1224.6923
$ws.Range("J116").Value = 2546.6667
$ws.Range("K116").Value = 1224.6923
$ws.Range("L116").Value = 2546.6667
$ws.Range("M116").Value = 1069.3077
$ws.Range("N116").Value = -7134.6667

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1642.1578
$ws.Range("I3").Value = 1224.6923
$ws.Range("J3").Value = 2546.6667
$ws.Range("K3").Value = 1224.6923
$ws.Range("L3").Value = 2546.6667
$ws.Range("M3").Value = -1110.6923
$ws.Range("N3").Value = -2774.6667

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 494.7143
$ws.Range("I22").Value = 300.75
$ws.Range("J22").Value = 753.3333
$ws.Range("K22").Value = 300.75
$ws.Range("L22").Value = 753.3333
$ws.Range("M22").Value = 49.25
$ws.Range("N22").Value = -1453.3333
$ws.Range("H99").Value = 2977.8
$ws.Range("I99").Value = 2517.3333
$ws.Range("J99").Value = 3354.5454
$ws.Range("K99").Value = 2517.3333
$ws.Range("L99").Value = 3354.5454
$ws.Range("M99").Value = -1019.3333
$ws.Range("N99").Value = -6350.5454
$ws.Range("H122").Value = 58824770
$ws.Range("I122").Value = 66667540
$ws.Range("K122").Value = 200002620
$ws.Range("M122").Value = -200000170
$ws.Range("H126").Value = 2977.8
$ws.Range("I126").Value = 2517.3333
$ws.Range("J126").Value = 3354.5454
$ws.Range("K126").Value = 7551.999899999999
$ws.Range("L126").Value = 10063.6362
$ws.Range("M126").Value = -5081.999899999999
$ws.Range("N126").Value = -15003.6362
$ws.Range("H132").Value = 1825.6364
$ws.Range("I132").Value = 1334.6154
$ws.Range("J132").Value = 2265.862
$ws.Range("K132").Value = 4003.8462
$ws.Range("L132").Value = 6797.586
$ws.Range("M132").Value = -1473.8462
$ws.Range("N132").Value = -11857.586
$ws.Range("H134").Value = 1947.975
$ws.Range("I134").Value = 1141.0555
$ws.Range("J134").Value = 2608.182
$ws.Range("K134").Value = 3423.1665
$ws.Range("L134").Value = 7824.545999999999
$ws.Range("M134").Value = -888.1664999999998
$ws.Range("N134").Value = -12894.546

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 958.7727
$ws.Range("I92").Value = 713.7143
$ws.Range("K92").Value = 2141.1429
$ws.Range("M92").Value = -893.1428999999998
$ws.Range("H113").Value = 37037932
$ws.Range("I113").Value = 637.5
$ws.Range("J113").Value = 66667770
$ws.Range("K113").Value = 1912.5
$ws.Range("L113").Value = 200003310
$ws.Range("M113").Value = 257.5
$ws.Range("N113").Value = -200007650

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 58.166668
$ws.Range("I2").Value = 58.42857
$ws.Range("J2").Value = 57.8
$ws.Range("K2").Value = 58.42857
$ws.Range("L2").Value = 57.8
$ws.Range("M2").Value = 54.57143
$ws.Range("N2").Value = -283.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 556107.25
$ws.Range("I46").Value = 500
$ws.Range("J46").Value = 1111714.5
$ws.Range("K46").Value = 500
$ws.Range("L46").Value = 1111714.5
$ws.Range("M46").Value = -312
$ws.Range("N46").Value = -1112090.5
$ws.Range("H61").Value = 4201.364
$ws.Range("I61").Value = 4651.875
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 4651.875
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -4449.875
$ws.Range("N61").Value = -3404
$ws.Range("H64").Value = 29750
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 29750
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 29750
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -30200
$ws.Range("H67").Value = 29750
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 29750
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 29750
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -31310
$ws.Range("H113").Value = 4201.364
$ws.Range("I113").Value = 4651.875
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 4651.875
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -2481.875
$ws.Range("N113").Value = -7340

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 17000
$ws.Range("J63").Value = 17000
$ws.Range("L63").Value = 17000
$ws.Range("N63").Value = -18248
$ws.Range("H66").Value = 17000
$ws.Range("J66").Value = 17000
$ws.Range("L66").Value = 51000
$ws.Range("N66").Value = -57240
$ws.Range("H122").Value = 1456.1724
$ws.Range("I122").Value = 1523
$ws.Range("K122").Value = 4569
$ws.Range("M122").Value = -2119
